# Adding different messages to the InputException for different input types
# -> appends two new rows of follow-up log data to the Follow_Up_Log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4: Subject "k" gets a follow-up entry (no follow-up date yet, notes "j")
$ws.Range("A4").Value = "k"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "2017-10-28"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "22:02:14.908319"
$ws.Range("D4").Value = "k"
$ws.Range("E4").Value = "j"
$ws.Range("F4").Value = "j"
$ws.Range("G4").Value = "j"

# New row 5: Subject "j" with a completed follow-up (date + time filled in, notes "n")
$ws.Range("A5").Value = "j"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "2017-10-28"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "22:04:29.674747"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "2017-10-28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "22:04:34.414847"
$ws.Range("F5").Value = "n"
$ws.Range("G5").Value = "j"
